$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ------------------------------------------------------------------
# Add 5 new rows (271-275) of "Bali" live-cam entries, mirroring the
# layout/style of the existing rows by cloning row 270's formatting
# first, then overwriting the cell contents.
# ------------------------------------------------------------------
for ($i = 271; $i -le 275; $i++) {
    $ws.Range("A270:G270").Copy($ws.Range("A" + $i + ":G" + $i))
}

# Row 271 - Dream Land Beach
$ws.Range("A271").Value = "LIVE, SEA, BEACH"
$ws.Range("B271").Value = "-8.79943111333464, 115.11765597255722"
$ws.Range("C271").Value = "Bali - Dream Land Beach, LIVE Webcam, BaliForum"
$ws.Range("D271").Value = "Bali"
$ws.Range("E271").Value = "Indonesia"
$ws.Range("F271").Value = "9yT659mJKR4"
$ws.Range("G271").Value = "Bali - Dream Land Beach, LIVE Webcam, BaliForum"

# Row 272 - Lucky Fish Lounge, Bingin Beach
$ws.Range("A272").Value = "LIVE, SEA, BEACH"
$ws.Range("B272").Value = "-8.80514003136993, 115.11360623416108"
$ws.Range("C272").Value = "Lucky Fish Lounge, Bingin Beach"
$ws.Range("D272").Value = "Bali"
$ws.Range("E272").Value = "Indonesia"
$ws.Range("F272").Value = "UUA2QTmTNaM"
$ws.Range("G272").Value = "Bali - Dream Land Beach, LIVE Webcam, BaliForum"

# Row 273 - Jimbaran Beach
$ws.Range("A273").Value = "LIVE, SEA, BEACH"
$ws.Range("B273").Value = "-8.780932571649334, 115.16380907469706"
$ws.Range("C273").Value = "Jimbaran Beach, Live Webcam. BaliForum & AKUSUKA Café"
$ws.Range("D273").Value = "Bali"
$ws.Range("E273").Value = "Indonesia"
$ws.Range("F273").Value = "mvVoilECpoY"
$ws.Range("G273").Value = "Bali - Dream Land Beach, LIVE Webcam, BaliForum"

# Row 274 - Jemeluk Amed
$ws.Range("A274").Value = "LIVE, SEA, BEACH"
$ws.Range("B274").Value = "-8.337822661077277, 115.66009343901486"
$ws.Range("C274").Value = "Jemeluk Amed, Live Webcam. BaliForum & See you again Café"
$ws.Range("D274").Value = "Bali"
$ws.Range("E274").Value = "Indonesia"
$ws.Range("F274").Value = "DQR8yF_FNwA"
$ws.Range("G274").Value = "Bali - Dream Land Beach, LIVE Webcam, BaliForum"

# Row 275 - Tegallalang rice fields
$ws.Range("A275").Value = "LIVE, MOUNTAIN, NATURE"
$ws.Range("B275").Value = "-8.433472622758597, 115.27892183636139"
$ws.Range("C275").Value = "Tegallalang rice fields, Bali. Online webcam. BaliForum & Rice terrace café"
$ws.Range("D275").Value = "Bali"
$ws.Range("E275").Value = "Indonesia"
$ws.Range("F275").Value = "yRPflX87vj0"
$ws.Range("G275").Value = "Bali - Dream Land Beach, LIVE Webcam, BaliForum"

# ------------------------------------------------------------------
# Hyperlink the new YouTube-link cells (G271:G275) to the BaliForum
# channel, keeping the existing title text as the visible cell text.
# (NB: don't read .Value back off a Range here - just reuse the
# literal we already wrote above - the host's Range.Value getter does
# not round-trip plain strings reliably.)
# ------------------------------------------------------------------
$baliChannel = "https://www.youtube.com/@BaliForumVideo"
$dreamLandTitle = "Bali - Dream Land Beach, LIVE Webcam, BaliForum"

$ws.Hyperlinks.Add($ws.Range("G271"), $baliChannel, "", "", $baliChannel) | Out-Null
$ws.Range("G271").Value = $dreamLandTitle
$ws.Range("G271").VerticalAlignment = -4108

$ws.Hyperlinks.Add($ws.Range("G272"), $baliChannel, "", "", $baliChannel) | Out-Null
$ws.Range("G272").Value = $dreamLandTitle
$ws.Range("G272").VerticalAlignment = -4108

$ws.Hyperlinks.Add($ws.Range("G273"), $baliChannel, "", "", $baliChannel) | Out-Null
$ws.Range("G273").Value = $dreamLandTitle
$ws.Range("G273").VerticalAlignment = -4108

$ws.Hyperlinks.Add($ws.Range("G274"), $baliChannel, "", "", $baliChannel) | Out-Null
$ws.Range("G274").Value = $dreamLandTitle
$ws.Range("G274").VerticalAlignment = -4108

$ws.Hyperlinks.Add($ws.Range("G275"), $baliChannel, "", "", $baliChannel) | Out-Null
$ws.Range("G275").Value = $dreamLandTitle
$ws.Range("G275").VerticalAlignment = -4108

# ------------------------------------------------------------------
# Update the selection bookkeeping to reflect the newly added rows
# (mirrors Excel's own "select the next empty row" behaviour after
# entering new data).
# ------------------------------------------------------------------
$ws.Range("A276").Select() | Out-Null
